# "added prob 20 in common part"
# Adds a new top-level problem summary (m0006, row 6) and a new set of
# "common part" solution steps (d0029 .. d0034, rows 79-84) underneath the
# existing d0027/d0028 pair (rows 77-78), fixes a typo in the d0028
# description, adds a missing formula cell at C44, and pushes the trailing
# single-cell x0001/y0001/z0001 marker rows down by one slot each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Push the x0001 / y0001 / z0001 marker rows down (bottom-up so we never
#     clobber a value before it has been copied onward). Row 81 will be
#     reused below for the new d0031 step, so its old "x0001" marker moves
#     to row 91 (etc.) ahead of time.
$ws.Range("A111").Value = "z0001"
$ws.Range("A101").Value = "y0001"
$ws.Range("A91").Value  = "x0001"

# --- Fix the typo in the existing d0028 description (이차대부등식 -> 이차부등식)
$ws.Range("B78").Value = '이차부등식이 항상 성립하도록 판별식에 대한 부등식을 세웁니다.'

# --- New formula cell completing row 44 (c0034 / sum identity row)
$ws.Range("C44").Value = '$\displaystyle\sum_{k=1}^{10} 2 a_{k}-\displaystyle\sum_{k=1}^{8} a_{k}=100$;'

# --- New "common part" problem 20 solution steps (rows 79-84)
$ws.Range("A79").Value = "d0029"
$ws.Range("C79").Value = '$0 \leq a \leq 6$;'
$ws.Range("B79").Value = '부등식을 풀어서 조건에 맞는 최댓값을 구합니다.'

$ws.Range("B80").Value = '부등식을 풀어서 조건에 맞는 최솟값을 구합니다.'
$ws.Range("A80").Value = "d0030"

$ws.Range("A81").Value = "d0031"
$ws.Range("B81").Value = '수치대입법으로 항등식의 미정계수값을 구합니다.'

$ws.Range("A82").Value = "d0032"
$ws.Range("C81").Value = '$b=1$;'
$ws.Range("B82").Value = '$x+1=t$ 로 치환해서 구간 $[1, 2]$에서의 함수 $f(t)$를 구합니다.'
$ws.Range("C82").Value = '$f(t)$;'

$ws.Range("A83").Value = "d0033"
$ws.Range("B83").Value = '함수가 미분가능하도록 미정계수값을 구합니다.'
$ws.Range("C83").Value = '$a=1$;'

$ws.Range("A84").Value = "d0034"
$ws.Range("B84").Value = '문제에 주어진 정적분을 계산합니다.'
$ws.Range("C84").Value = '$ \displaystyle\int_{1}^{2} f(x) d x$;'

# --- New top-level problem summary row (m0006, row 6)
$ws.Range("A6").Value = "m0006"
$ws.Range("B6").Value = '지금까지 구해진 값들을 이용해서 문제에서 요구하는 값을 구합니다.'
$ws.Range("C6").Value = '$60 \times \displaystyle\int_{1}^{2} f(x) d x$;'

# --- Update the view to reflect the newly selected/visible cell
$ws.Range("A86").Select()
